$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 725.8889
$ws.Cells.Item(12, 9).Value = 869
$ws.Cells.Item(12, 11).Value = 869
$ws.Cells.Item(12, 13).Value = -699
$ws.Cells.Item(33, 8).Value = 330.4516
$ws.Cells.Item(33, 9).Value = 341.8889
$ws.Cells.Item(33, 10).Value = 253.25
$ws.Cells.Item(33, 11).Value = 341.8889
$ws.Cells.Item(33, 12).Value = 253.25
$ws.Cells.Item(33, 13).Value = -112.8889
$ws.Cells.Item(33, 14).Value = -711.25
$ws.Cells.Item(99, 8).Value = 1479.3334
$ws.Cells.Item(99, 9).Value = 1620
$ws.Cells.Item(99, 10).Value = 776
$ws.Cells.Item(99, 11).Value = 4860
$ws.Cells.Item(99, 12).Value = 2328
$ws.Cells.Item(99, 13).Value = -3362
$ws.Cells.Item(99, 14).Value = -5324
$ws.Cells.Item(112, 8).Value = 3265.5
$ws.Cells.Item(112, 10).Value = 3390
$ws.Cells.Item(112, 12).Value = 10170
$ws.Cells.Item(112, 14).Value = -12386
$ws.Cells.Item(129, 8).Value = 3044.5557
$ws.Cells.Item(129, 9).Value = 898
$ws.Cells.Item(129, 11).Value = 2694
$ws.Cells.Item(129, 13).Value = 2306
$ws.Cells.Item(137, 8).Value = 1529.7727
$ws.Cells.Item(137, 9).Value = 1328.1538
$ws.Cells.Item(137, 10).Value = 1821
$ws.Cells.Item(137, 11).Value = 3984.4614
$ws.Cells.Item(137, 12).Value = 5463
$ws.Cells.Item(137, 13).Value = -1434.4614
$ws.Cells.Item(137, 14).Value = -10563
$ws.Cells.Item(138, 8).Value = 2700.7827
$ws.Cells.Item(138, 9).Value = 1707.8667
$ws.Cells.Item(138, 11).Value = 5123.6001
$ws.Cells.Item(138, 13).Value = 16.39990000000034

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(61, 8).Value = 0
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 11).Value = 0
$ws.Cells.Item(61, 13).ClearContents()
$ws.Cells.Item(97, 8).Value = 1763.6666
$ws.Cells.Item(97, 9).Value = 621.3333
$ws.Cells.Item(97, 11).Value = 621.3333
$ws.Cells.Item(97, 13).Value = -125.3333
$ws.Cells.Item(136, 8).Value = 0
$ws.Cells.Item(136, 9).Value = 0
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 13).ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 222.5
$ws.Cells.Item(80, 8).Value = 351.7857
$ws.Cells.Item(80, 9).Value = 548.6
$ws.Cells.Item(80, 11).Value = 548.6
$ws.Cells.Item(80, 13).Value = 449.4
$ws.Cells.Item(83, 8).Value = 351.7857
$ws.Cells.Item(83, 9).Value = 548.6
$ws.Cells.Item(83, 11).Value = 2743
$ws.Cells.Item(83, 13).Value = 2249
$ws.Cells.Item(134, 8).Value = 1828
$ws.Cells.Item(134, 10).Value = 0
$ws.Cells.Item(134, 12).Value = 0
$ws.Cells.Item(134, 14).ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value = 2018.8572
$ws.Cells.Item(58, 9).Value = 1426.4
$ws.Cells.Item(58, 11).Value = 1426.4
$ws.Cells.Item(58, 13).Value = -1223.4
$ws.Cells.Item(107, 8).Value = 990.9
$ws.Cells.Item(107, 9).Value = 772.8570999999999
$ws.Cells.Item(107, 11).Value = 772.8570999999999
$ws.Cells.Item(107, 13).Value = 1147.1429
$ws.Cells.Item(134, 8).Value = 1528.7307
$ws.Cells.Item(134, 9).Value = 1138.7142
$ws.Cells.Item(134, 10).Value = 3166.8
$ws.Cells.Item(134, 11).Value = 3416.1426
$ws.Cells.Item(134, 12).Value = 9500.400000000001
$ws.Cells.Item(134, 13).Value = -881.1425999999997
$ws.Cells.Item(134, 14).Value = -14570.4
$ws.Cells.Item(136, 8).Value = 2018.8572
$ws.Cells.Item(136, 9).Value = 1426.4
$ws.Cells.Item(136, 11).Value = 4279.200000000001
$ws.Cells.Item(136, 13).Value = -1729.200000000001
$ws.Cells.Item(141, 8).Value = 150000
$ws.Cells.Item(141, 10).Value = 150000
$ws.Cells.Item(141, 12).Value = 150000
$ws.Cells.Item(141, 14).Value = -160360

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(38, 8).Value = 87.55556
$ws.Cells.Item(38, 10).Value = 112.666664
$ws.Cells.Item(38, 12).Value = 337.999992
$ws.Cells.Item(38, 14).Value = -1031.999992
$ws.Cells.Item(68, 8).Value = 0
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 13).ClearContents()
$ws.Cells.Item(71, 8).Value = 0
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 13).ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 2677.8333
$ws.Cells.Item(80, 9).Value = 2548.2
$ws.Cells.Item(80, 10).Value = 2770.4285
$ws.Cells.Item(80, 11).Value = 2548.2
$ws.Cells.Item(80, 12).Value = 2770.4285
$ws.Cells.Item(80, 13).Value = -1550.2
$ws.Cells.Item(80, 14).Value = -4766.4285
$ws.Cells.Item(83, 8).Value = 2677.8333
$ws.Cells.Item(83, 9).Value = 2548.2
$ws.Cells.Item(83, 10).Value = 2770.4285
$ws.Cells.Item(83, 11).Value = 12741
$ws.Cells.Item(83, 12).Value = 13852.1425
$ws.Cells.Item(83, 13).Value = -7749
$ws.Cells.Item(83, 14).Value = -23836.1425
$ws.Cells.Item(107, 8).Value = 13748
$ws.Cells.Item(107, 9).Value = 1600
$ws.Cells.Item(107, 10).Value = 17797.334
$ws.Cells.Item(107, 11).Value = 1600
$ws.Cells.Item(107, 12).Value = 17797.334
$ws.Cells.Item(107, 13).Value = 320
$ws.Cells.Item(107, 14).Value = -21637.334
$ws.Cells.Item(113, 8).Value = 1590.8
$ws.Cells.Item(113, 9).Value = 1590.8
$ws.Cells.Item(113, 11).Value = 1590.8
$ws.Cells.Item(113, 13).Value = 579.2
$ws.Cells.Item(126, 8).Value = 3418
$ws.Cells.Item(126, 9).Value = 3647.5
$ws.Cells.Item(126, 10).Value = 2500
$ws.Cells.Item(126, 11).Value = 10942.5
$ws.Cells.Item(126, 12).Value = 7500
$ws.Cells.Item(126, 13).Value = -8472.5
$ws.Cells.Item(126, 14).Value = -12440

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1726.5264
$ws.Cells.Item(46, 10).Value = 1940.4
$ws.Cells.Item(46, 12).Value = 1940.4
$ws.Cells.Item(46, 14).Value = -2316.4
$ws.Cells.Item(55, 8).Value = 554.75
$ws.Cells.Item(55, 10).Value = 691.6667
$ws.Cells.Item(55, 12).Value = 691.6667
$ws.Cells.Item(55, 14).Value = -1037.6667
$ws.Cells.Item(74, 8).Value = 39999.668
$ws.Cells.Item(74, 10).Value = 39999.668
$ws.Cells.Item(74, 12).Value = 39999.668
$ws.Cells.Item(74, 14).Value = -41995.668
$ws.Cells.Item(77, 8).Value = 39999.668
$ws.Cells.Item(77, 10).Value = 39999.668
$ws.Cells.Item(77, 12).Value = 119999.004
$ws.Cells.Item(77, 14).Value = -129983.004
$ws.Cells.Item(82, 8).Value = 0
$ws.Cells.Item(82, 10).Value = 0
$ws.Cells.Item(82, 12).Value = 0
$ws.Cells.Item(82, 14).ClearContents()
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 12).Value = 0
$ws.Cells.Item(85, 14).ClearContents()
$ws.Cells.Item(122, 8).Value = 6543
$ws.Cells.Item(122, 9).Value = 5621.625
$ws.Cells.Item(122, 10).Value = 9000
$ws.Cells.Item(122, 11).Value = 16864.875
$ws.Cells.Item(122, 12).Value = 27000
$ws.Cells.Item(122, 13).Value = -14414.875
$ws.Cells.Item(122, 14).Value = -31900
$ws.Cells.Item(132, 8).Value = 1628.2307
$ws.Cells.Item(132, 9).Value = 1531.3636
$ws.Cells.Item(132, 11).Value = 4594.0908
$ws.Cells.Item(132, 13).Value = -2064.0908

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(82, 8).Value = 0
$ws.Cells.Item(82, 9).Value = 0
$ws.Cells.Item(82, 11).Value = 0
$ws.Cells.Item(82, 13).ClearContents()
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 11).Value = 0
$ws.Cells.Item(85, 13).ClearContents()
$ws.Cells.Item(122, 8).Value = 1433.7142
$ws.Cells.Item(122, 9).Value = 1561.3636
$ws.Cells.Item(122, 10).Value = 965.6667
$ws.Cells.Item(122, 11).Value = 4684.0908
$ws.Cells.Item(122, 12).Value = 2897.0001
$ws.Cells.Item(122, 13).Value = -2234.0908
$ws.Cells.Item(122, 14).Value = -7797.0001
$ws.Cells.Item(132, 8).Value = 1534
$ws.Cells.Item(132, 9).Value = 1563.7222
$ws.Cells.Item(132, 10).Value = 999
$ws.Cells.Item(132, 11).Value = 4691.1666
$ws.Cells.Item(132, 12).Value = 2997
$ws.Cells.Item(132, 13).Value = -2161.1666
$ws.Cells.Item(132, 14).Value = -8057
